$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.607.91'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '3.765.29'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'594.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = "'167.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("D7").Value = '3.763.49'
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("D11").Value = "'6.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.60%  '
$ws.Range("D14").Value = "'36.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").Value = '4.397.25'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '3.753.28'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '68.595.00'
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("D18").Value = "'17.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.95%  '
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = "'7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.87%  '
$ws.Range("D21").Value = "'10.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("D22").Value = "'466.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").Value = "'0.700"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = "'0.0000149"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = "'84.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("E26").Value = '  -2.79%  '
$ws.Range("D27").Value = "'11.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").Value = "'10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.70%  '
$ws.Range("D30").Value = '3.912.14'
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -4.56%  '
$ws.Range("D32").Value = "'7.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.47%  '
$ws.Range("D33").Value = "'30.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("D34").Value = "'2.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("D35").Value = "'9.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = '3.724.80'
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("E38").Value = '  -3.55%  '
$ws.Range("D39").Value = "'3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.54%  '
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("D42").Value = "'5.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = "'44.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.66%  '
$ws.Range("E46").Value = '  -3.39%  '
$ws.Range("D47").Value = "'46.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D49").Value = "'8.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = "'391.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = "'145.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.68%  '
